$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.15576752824467
$ws.Range("C2").Value = 9.920862905568981
$ws.Range("E2").Value = 11.72031136012316
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.631885520300537
$ws.Range("L2").Value = 9.712552853149091
$ws.Range("N2").Value = 17.35157256621244
$ws.Range("O2").Value = 22.20859199889635
$ws.Range("B3").Value = 15.67285348670882
$ws.Range("C3").Value = 9.771982714641043
$ws.Range("E3").Value = 11.75636112744111
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.634048106775754
$ws.Range("L3").Value = 9.68824903374883
$ws.Range("N3").Value = 17.4032326530195
$ws.Range("O3").Value = 22.27284906520507
$ws.Range("B4").Value = 15.37117320656519
$ws.Range("C4").Value = 9.678735950586361
$ws.Range("E4").Value = 11.78056535155947
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.635446395201351
$ws.Range("L4").Value = 9.675024872996669
$ws.Range("N4").Value = 17.4368938538647
$ws.Range("O4").Value = 22.31849814783111
$ws.Range("B5").Value = 15.24712734068469
$ws.Range("C5").Value = 9.640304344299116
$ws.Range("E5").Value = 11.79094878139155
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.636033981710043
$ws.Range("L5").Value = 9.670066858221071
$ws.Range("N5").Value = 17.45110024810711
$ws.Range("O5").Value = 22.33865207569072
$ws.Range("B6").Value = 15.22646849594758
$ws.Range("C6").Value = 9.633897508463903
$ws.Range("E6").Value = 11.79270433787227
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.636132625102537
$ws.Range("L6").Value = 9.669269716163917
$ws.Range("N6").Value = 17.45348878195582
$ws.Range("O6").Value = 22.34209214118792
$ws.Range("B7").Value = 15.36950450765208
$ws.Range("C7").Value = 9.678219362879688
$ws.Range("E7").Value = 11.78070328113402
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.635454247561832
$ws.Range("L7").Value = 9.674956257825491
$ws.Range("N7").Value = 17.43708346435701
$ws.Range("O7").Value = 22.31876367708622
$ws.Range("B8").Value = 15.99044218978649
$ws.Range("C8").Value = 9.869926333877933
$ws.Range("E8").Value = 11.73231156404388
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.632616592245702
$ws.Range("L8").Value = 9.703822789181855
$ws.Range("N8").Value = 17.36898251837905
$ws.Range("O8").Value = 22.22945843270647
$ws.Range("B9").Value = 17.15936900668395
$ws.Range("C9").Value = 10.23011502069406
$ws.Range("E9").Value = 11.65385309216506
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.627608358593055
$ws.Range("L9").Value = 9.773722445019393
$ws.Range("N9").Value = 17.25080183104059
$ws.Range("O9").Value = 22.10374819900287
$ws.Range("B10").Value = 17.97911894639157
$ws.Range("C10").Value = 10.48361041841379
$ws.Range("E10").Value = 11.6062527306
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.624264362688184
$ws.Range("L10").Value = 9.832912350592995
$ws.Range("N10").Value = 17.17328456520564
$ws.Range("O10").Value = 22.04185015441163
$ws.Range("B11").Value = 18.3418640428329
$ws.Range("C11").Value = 10.59620989246014
$ws.Range("E11").Value = 11.58678330432595
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.6228151793917
$ws.Range("L11").Value = 9.861476087574347
$ws.Range("N11").Value = 17.14002985199509
$ws.Range("O11").Value = 22.0203682983644
$ws.Range("B12").Value = 18.47764137718035
$ws.Range("C12").Value = 10.6384336906295
$ws.Range("E12").Value = 11.57972517747155
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.622276708562183
$ws.Range("L12").Value = 9.872521969165886
$ws.Range("N12").Value = 17.12772510623854
$ws.Range("O12").Value = 22.01319802132613
$ws.Range("B13").Value = 18.44847170480014
$ws.Range("C13").Value = 10.62935888259633
$ws.Range("E13").Value = 11.58123127626395
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.622392220474278
$ws.Range("L13").Value = 9.870132938429874
$ws.Range("N13").Value = 17.13036235122462
$ws.Range("O13").Value = 22.01469930986444
$ws.Range("B14").Value = 18.35306698854416
$ws.Range("C14").Value = 10.59969212657064
$ws.Range("E14").Value = 11.58619632208673
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.622770672875829
$ws.Range("L14").Value = 9.86238027726151
$ws.Range("N14").Value = 17.1390117632535
$ws.Range("O14").Value = 22.01975904297067
$ws.Range("B15").Value = 18.29441876735857
$ws.Range("C15").Value = 10.58146560691394
$ws.Range("E15").Value = 11.58927852857412
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.62300382615839
$ws.Range("L15").Value = 9.857661237972652
$ws.Range("N15").Value = 17.14434727577392
$ws.Range("O15").Value = 22.02298399689547
$ws.Range("B16").Value = 17.95519721768602
$ws.Range("C16").Value = 10.47619498619731
$ws.Range("E16").Value = 11.60756909314955
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.624360514893787
$ws.Range("L16").Value = 9.831078086915003
$ws.Range("N16").Value = 17.1754981906053
$ws.Range("O16").Value = 22.04338877230596
$ws.Range("B17").Value = 17.74439972444802
$ws.Range("C17").Value = 10.41090120647135
$ws.Range("E17").Value = 11.61934950958955
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.625211207305802
$ws.Range("L17").Value = 9.815185448975198
$ws.Range("N17").Value = 17.1951221336526
$ws.Range("O17").Value = 22.05761962391168
$ws.Range("B18").Value = 17.62220458692905
$ws.Range("C18").Value = 10.37309161794066
$ws.Range("E18").Value = 11.626330857952
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.625707284797877
$ws.Range("L18").Value = 9.806198982674061
$ws.Range("N18").Value = 17.20659837545528
$ws.Range("O18").Value = 22.06643296385602
$ws.Range("B19").Value = 17.58067226922737
$ws.Range("C19").Value = 10.36024702509129
$ws.Range("E19").Value = 11.62872991478342
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.625876414488105
$ws.Range("L19").Value = 9.803183046533746
$ws.Range("N19").Value = 17.21051652810857
$ws.Range("O19").Value = 22.06952474729225
$ws.Range("B20").Value = 17.76693874412538
$ws.Range("C20").Value = 10.41787834463449
$ws.Range("E20").Value = 11.61807418577883
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.625119948145782
$ws.Range("L20").Value = 9.816861294515819
$ws.Range("N20").Value = 17.19301356860947
$ws.Range("O20").Value = 22.05603968239163
$ws.Range("B21").Value = 18.3811336513058
$ws.Range("C21").Value = 10.60841742901983
$ws.Range("E21").Value = 11.58472942839127
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.622659233016345
$ws.Range("L21").Value = 9.864651245858383
$ws.Range("N21").Value = 17.13646340795262
$ws.Range("O21").Value = 22.01824666928507
$ws.Range("B22").Value = 18.77325000138203
$ws.Range("C22").Value = 10.73051525062641
$ws.Range("E22").Value = 11.564770162084
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.621111045617647
$ws.Range("L22").Value = 9.89721888954714
$ws.Range("N22").Value = 17.10118352658284
$ws.Range("O22").Value = 21.99916991461622
$ws.Range("B23").Value = 18.56485863622052
$ws.Range("C23").Value = 10.66557951426973
$ws.Range("E23").Value = 11.57525489945887
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.62193186668166
$ws.Range("L23").Value = 9.879716952074352
$ws.Range("N23").Value = 17.11985966567386
$ws.Range("O23").Value = 22.00883568164874
$ws.Range("B24").Value = 17.75675197356493
$ws.Range("C24").Value = 10.41472482374845
$ws.Range("E24").Value = 11.61865010984826
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.625161184594012
$ws.Range("L24").Value = 9.816103175738398
$ws.Range("N24").Value = 17.19396624594844
$ws.Range("O24").Value = 22.05675200613475
$ws.Range("B25").Value = 16.84938040110588
$ws.Range("C25").Value = 10.13452281185079
$ws.Range("E25").Value = 11.67331605391826
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.628904030105664
$ws.Range("L25").Value = 9.753416010170563
$ws.Range("N25").Value = 17.28113387325007
$ws.Range("O25").Value = 22.13242937661439
